$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q3" right after "总计" (i.e. before
#    the current "2022-Q2" sheet), pushing every later sheet one slot right.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet   = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# -- Header row ----------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# -- Data rows ---------------------------------------------------------
# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "007012"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "湘财长顺混合A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "3.12"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "94.24"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "6.35"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.1981"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 6

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "007013"
$newSheet.Range("B3").Style = "Normal"
$newSheet.Range("C3").Value = "湘财长顺混合C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "1.29"
$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "94.24"
$newSheet.Range("E3").Style = "Normal"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "6.35"
$newSheet.Range("F3").Style = "Normal"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0819"
$newSheet.Range("G3").Style = "Normal"
$newSheet.Range("H3").Value = 6

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").Value = "008128"
$newSheet.Range("B4").Style = "Normal"
$newSheet.Range("C4").Value = "湘财长源股票A"
$newSheet.Range("D4").NumberFormat = "@"
$newSheet.Range("D4").Value = "1.32"
$newSheet.Range("D4").Style = "Normal"
$newSheet.Range("E4").NumberFormat = "@"
$newSheet.Range("E4").Value = "94.07"
$newSheet.Range("E4").Style = "Normal"
$newSheet.Range("F4").NumberFormat = "@"
$newSheet.Range("F4").Value = "5.17"
$newSheet.Range("F4").Style = "Normal"
$newSheet.Range("G4").NumberFormat = "@"
$newSheet.Range("G4").Value = "0.0682"
$newSheet.Range("G4").Style = "Normal"
$newSheet.Range("H4").Value = 10

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").NumberFormat = "@"
$newSheet.Range("B5").Value = "008129"
$newSheet.Range("B5").Style = "Normal"
$newSheet.Range("C5").Value = "湘财长源股票C"
$newSheet.Range("D5").NumberFormat = "@"
$newSheet.Range("D5").Value = "0.53"
$newSheet.Range("D5").Style = "Normal"
$newSheet.Range("E5").NumberFormat = "@"
$newSheet.Range("E5").Value = "94.07"
$newSheet.Range("E5").Style = "Normal"
$newSheet.Range("F5").NumberFormat = "@"
$newSheet.Range("F5").Value = "5.17"
$newSheet.Range("F5").Style = "Normal"
$newSheet.Range("G5").NumberFormat = "@"
$newSheet.Range("G5").Value = "0.0274"
$newSheet.Range("G5").Style = "Normal"
$newSheet.Range("H5").Value = 10

# -- Formatting (bold / thin border / centred+top-aligned), matching the
#    look of the header + index column on the sibling quarter sheets -------
$hdr = $newSheet.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$idxCol = $newSheet.Range("A2:A5")
$idxCol.Font.Bold = $true
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
$idxCol.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the 2022-Q3 row at the top of
#    the data (row 2) and shift the remaining rows down by one, renumbering
#    the leading index column (A) 0..5.
# ---------------------------------------------------------------------------
$ws = $totalSheet

# Shift existing rows 6->7, 5->6, 4->5, 3->4, 2->3 (bottom-up so values
# aren't clobbered before being read), then write the brand-new row 2.
$ws.Range("B7").Value = "2021-Q1"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 0.52

$ws.Range("B6").Value = "2021-Q2"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 0.54

$ws.Range("B5").Value = "2021-Q3"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 0.16

$ws.Range("B4").Value = "2021-Q4"
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 0.18

$ws.Range("B3").Value = "2022-Q2"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0

$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 0.38

# Renumber the index column.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# Extend the "A" index-column style (s="2") down to the new last row (7) -
# copy from the row above, which already carries the correct style. This
# must be the final thing touching A6/A7 so the pasted format isn't clobbered
# by a later write.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
